# Work on jump mechanic
# Append a new day-log entry (row 4) to the time-tracking sheet:
#   Wednesday | 23-Apr-2025 | 11:00 -> 11:30 | 0:30 | "Started on jump"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day
$ws.Range("A4").Value = "Wednesday"

# Date (format matches the existing "d-mmm" style used for this new row)
$ws.Range("B4").Value = 45770
$ws.Range("B4").NumberFormat = "d-mmm"

# From / Until / Time spent - reuse the same time format as the row above
$ws.Range("C4").Value = 0.45833333333333331
$ws.Range("C4").NumberFormat = $ws.Range("C3").NumberFormat

$ws.Range("D4").Value = 0.47916666666666669
$ws.Range("D4").NumberFormat = $ws.Range("D3").NumberFormat

$ws.Range("E4").Value = 0.020833333333333332
$ws.Range("E4").NumberFormat = $ws.Range("E3").NumberFormat

# Realised
$ws.Range("F4").Value = "Started on jump"

# Move the selection past the newly added row, like the original author's sheet
$ws.Range("F5").Select() | Out-Null
